# This deck ships two theme parts:
#   ppt/theme/theme1.xml -> used by the (only) Slide Master, originally the
#                            "Integral" colour scheme
#   ppt/theme/theme2.xml -> used by the Notes Master, originally the
#                            stock "Office Theme" colour scheme
#
# The target revision re-colours the design applied to the slides (theme1)
# with the standard Office palette. Re-apply the 12 theme colours (dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink) on the Slide Master's theme so that
# ppt/theme/theme1.xml ends up holding the "Office" colour values that used
# to live in theme2.xml.

$p = $ppt.ActivePresentation

$design = $p.Designs.Item(1)
$colorScheme = $design.SlideMaster.Theme.ThemeColorScheme

# index -> RGB (as 0xBBGGRR long, matching the ColorFormat.RGB convention)
$colorScheme.Item(1).RGB  = 0        # dk1      000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink 954F72
